$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'The article discusses the findings of the Second National Climate Change Study, completed by Singapore''s MSS Centre for Climate Research. It highlights some of the key points as the complexity of generating climate scenarios and challenges in applying them, as well as the benefits of better regional cooperation and data sharing.'
$ws.Range("D3").Value = 'The UK''s climate action foreign aid reached an all-time high of £3bn in the latest year. However, the amount was controversial due to changes in how the UK calculates climate aid for developing countries. The article also reveals that the UK is on track to meet its five-year goal of £11.6bn in aid to developing countries by 2026, even though it continues to cut overall aid budgets. The article is a latest analysis of the UK''s international climate finance contributions since 2011.'
$ws.Range("D4").Value = 'The article reveals that N-tier suppliers (Tier 2-4) are a growing concern due to their suppliers being the root cause of significant supply chain disruptions. 
The lack of visibility into suppliers below Tier 1 is attributed to operational silos, where each level of supply chain management has its own separate systems and metrics, which prevent cohesive view and information gaps. 
Organizations can address this by implementing a more unified and collaborative approach, combining data from various systems, and establishing comprehensive risk assessments and proactive'
$ws.Range("D5").Value = 'The article highlights that the last decade has brought significant disruptions to supply chains, with rising risks involving bankruptcy, force majeure, and ESG.
The article explores how procurement leaders have reshaped their risk frameworks and investment strategies to address these massive disruptions, prioritizing supplier accountability and sustainability concerns.'
$ws.Range("D6").Value = 'Climate negotiators are trying to make progress and find solutions before COP30 while facing challenges such as funding for climate action, multilateralism, and geopolitical tensions. The key outcomes and disputes at the UNFCCC''s 62nd biannual session in Bonn include discussions around adaptation, the global stocktake, and just transition, with some progress on just-transition talks.'
$ws.Range("D7").Value = 'The Swiss government''s recent decision to pause further climate action plans is seen as a regression, lacking long-term vision for more sustainable practices and targets. The pause may lead to a delay in crucial emissions reductions and may result in missed opportunities for the country to lead broader climate action globally.'
$ws.Range("D8").Value = 'The article discusses strategies to remove carbon dioxide from the atmosphere, which can be done through land-based or ocean-based techniques. It highlights insights about the Read importance of understanding the benefits and risks of different CDR options, as well as the readiness, scalability and unknown factors of their impact on effectiveness, cost and potential impacts on people and ecosystems.'
$ws.Range("D9").Value = 'The article discusses how Chemical Recycling Europe partnered with sustainability firms to establish a standardized approach to life cycle assessments. This partnership helped develop a guidelines for assessing the environmental impacts of chemical recycling processes, promoting transparency and credible evaluations. The initiative helps industry to align methodologies, make assessments more consistent and support better decision-making and overall sustainability within the chemical recycling sector.'
$ws.Range("D10").Value = 'The article highlights how Endress+Hauser, a leader in industrial processes, partnered with Sphera to address customer demands for carbon footprints. Sphera helped experts provided guidance to develop a life cycle assessment model and software to help the client meet the challenge of providing accurate and reliable carbon footprint data to various product lines. This partnership streamlined the smar process, and Endress+Hauser is now equipped with a flexible and efficient tool for carbon footprint calculations, which sets a solid foundation for sustainability initiatives and'
$ws.Range("D11").Value = 'Under the new Labour government, the UK''s climate advisors now have more confidence that the nation can meet its emission targets. The report finds progress on important climate issues such also acknowledges areas of risk. The report also suggests lowering electricity prices to support industry and transportation electrification.'
$ws.Range("D12").Value = 'The article provides a concise summary of key climate change and environmental sustainability news from the week. It notes the developments like heat domes and wildfires in North America and Europe, as well as rising energy demand and power prices in Japan. It also discusses the budget growth at the Bonn climate talks and the agreement by many countries to increase the UN Framework Convention on Climate Change''s budget.'
$ws.Range("D13").Value = 'Alaska Ecosystem Status Reports compile and share results of various research projects related to the health of the ecosystem across the Gulf of Alaska, Bering Sea, and Aleutian Islands. This initiative aims to enhance our understanding of ecosystem health and improve sustainable management of fisheries.'
$ws.Range("D14").Value = 'The article talks about how companies can use software solutions to achieve sustainability goals, manage product stewardship, and ensure compliance with safety and regulatory standards. The software helps customers by providing expert understanding of industries and accurate data for managing chemical and material hazards.'
$ws.Range("D15").Value = 'Here''s a brief summary: The article highlights the importance of sustainability in Singapore''s energy sector, and how students can make a difference by pursuing energy industry scholarships. These scholarships are also open opportunities for students with engineering skills to make impactful contributions and have a positive impact on the energy sector while enjoying meaningful work and stability. The article also mentions scholarship eligibility and past scholars, focusing on the benefits of these opportunities.'
$ws.Range("D16").Value = 'The article discusses how Environmental Health and Safety (EHS) programs have evolved to encompass broader topics including environment, health, sustainability, and social responsibility. The focus is on managing and reporting on risks and hazards related to these topics, which have become critical for organizations.'
$ws.Range("D17").Value = 'Wolters Kluper is making waves by committing to sustainability and expanding its ESG offerings to create a more efficient and resilient future. The company aims to address critical global issues and drive positive change through data and intelligence solutions.'
$ws.Range("D18").Value = 'The article invites English-speaking climate experts from the global south to fill out a form to add their details to a database. It is a aim to include a diverse range of experts and perspectives in climate reporting. The term "global south" refers to lower-income countries in regions like Africa, Asia, and Latin America and is used to describe developing nations that have different histories, cultures, and economies. The term is not universally accepted and can cause confusion.'
$ws.Range("D19").Value = 'The key insights from this article are that volcanic eruptions have a substantial impact on climate change and climate models. The unpredictable and sporadic nature of volcanic eruptions poses challenges to current climate projections. The study found that including volcanic forcing in climate projections can slightly delay warming limits, but it also impacts the rate of warming and cooling over decades. The findings suggest that volcanic activity should be considered in climate projections to accurately reflect the complexities of climate change.'
$ws.Range("D20").Value = '1. The article reveals that Singapore had its hottest year on record in 2015, surpassing its average temperature by 0.8 degrees Celsius. It also highlights that without urgent action to cut emissions, the world will continue to see increasing temperatures and more extreme weather events. 
2. Singapore''s weather mobile app, �Weather@SG, launched today, provides users with access to relevant and useful weather and climate information.'
$ws.Range("D21").Value = 'The article reveals that 2017 was the warmest year on record, not influenced by El Niño, due to long-term temperature rise from global warming and urban development. The El Niño Southern Oscillation (ENSO), which influences temperatures and rainfall, was neutral in 2017, except for November and December. This shows that climate change in Singapore is primarily driven by factors beyond natural weather phenomena.'
$ws.Range("D22").Value = 'The article discusses the increasing concern of ocean acidification and its impacts on marine life, coastal areas, and the Great Lakes. It highlights the partnership programs as NOAA to help society adapt and respond to these growing threats from acidification on the ocean.'
$ws.Range("D23").Value = 'The article discusses how sustainability and safety management software solutions can help businesses meet goals, with a focus on product stewardship. The software helps manage hazardous materials and compliance, and provides accurate data for various industries, all critical for safe product manufacturing and proper regulations adherence.'
$ws.Range("D24").Value = 'The key insights from this article are that the world is set to witness an unprecedented rise in energy demand, emphasizing the urgent need for low-carbon energy solutions. The article also highlights the importance of energy transition to meet sustainability goals and the potential of carbon removal technologies to address the climate challenge.'
$ws.Range("D25").Value = 'The ESG article shares that the upcoming US presidential election will decide between incumbent Donald Trump and Joe Biden. The outcome will significantly impact the important policies and social issues.'
$ws.Range("D26").Value = 'The article highlights the key findings of Phase 1 of Singapore''s Second National Climate Change Study. It shows that the results are consistent with the IPCC 5th Assessment Report, which projects sea-level rise globally, along with higher temperatures and extreme rainfall for the Singapore region. The study aims to provide detailed projections for the climate changes that are of interest to Singapore, and is aimed to be completed soon.'
$ws.Range("D27").Value = 'The 2022 climate data reveals that last year was extremely hot and reached new temperature records. This underscores the Somewhat of a concern, these records come at a time when many countries are still trying to recover from the immediate impacts of the pandemic. This is also emphasizes the essential need for climate mitigation measures.'
$ws.Range("D28").Value = 'The article explores how ESG (Environmental, Social, and Governance) factors are crucial for investors to understand and incorporate into their strategies, and how search tools can help surface relevant content for better investment decisions. The piece emphasizes the particular attention to how new insights and concepts are driving the ESG space forward and creating opportunities for investors and companies.'
$ws.Range("D29").Value = 'The article highlights the urgent need for scientists to study glaciers and ice caps in the Arctic. The unique gases trapped within glacial ice provide historical information to scientists about atmospheric CO2 levels. This information is crucial for understanding climate change impacts.'
$ws.Range("D30").Value = 'The article discusses how the energy crisis has been navigated by the UK''s energy regulator. The CEO, Jonathan Brearley, has plans to transition the UK to clean energy by 2030, focusing on improving consumer protection against high energy prices. The CEO also notes the carbon pricing changes and a transition from gas heaters to clean energy usage.'
$ws.Range("D31").Value = 'The webinar covered insights from COP29, including the finalization of the loss and damage facility, the phase out of coal demand, and the emerging global emissions threshold.  The important will in addressing climate change and global warming, and addressing the issues of developing nations.'
$ws.Range("D32").Value = 'Yes, the article mentions that they have noticed changes in the permafrost landscape and biogeochemical research areas due to climate change impacts, such as melting permafrost and rising greenhouse gas emissions. They also discuss the particular research on climate change and the potential impact on vulnerable areas.'
